$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reshape the sheet: drop the title row and the blank rows, so the
#     table data starts at row 1 instead of row 5 -------------------------
$ws.Rows("1").Delete()
$ws.Range("1:3").Delete()

# --- Insert a new column (D) between "IP" and "LAN" for the new
#     "reseau" column -------------------------------------------------
$ws.Columns("D").Insert()

# --- Resize the table to include the new column and match the new
#     data extent ------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F6"))

# Re-assert the header labels so the table column names stay in sync
# with the worksheet header cells
$ws.Range("D1").Value2 = "réseau"
$ws.Range("E1").Value2 = "LAN"
$ws.Range("F1").Value2 = "Commentaire"

# --- Update the IP / LAN text that lost the OPT1/OPT2 entries -----------
$ws.Range("C2").Value2 = "WAN : 192.168.1.99`nLAN : 192.168.10.1`nWEB : 192.168.99.1"
$ws.Range("E2").Value2 = "lan :10`nOPT2: 99"

# --- Fill in the new "reseau" column values ------------------------------
$ws.Range("D2").Value2 = "WAN : 192.168.1.0/24`nLAN : 192.168.10.0/24`nWEB : 192.168.99.0/24"
$ws.Range("D3").Value2 = "192.168.10.0/24"
$ws.Range("D4").Value2 = "192.168.99.0/24"
$ws.Range("D5").Value2 = "192.168.1.0/24"
$ws.Range("D6").Value2 = "192.168.10.0/24"

# --- New row 6 : srvmail / server-AD duplicate ---------------------------
$ws.Range("A6").Value2 = "srvmail"
$ws.Range("C6").Value2 = "192.168.10.5"
$ws.Range("E6").Value2 = 10
$ws.Range("F6").Value2 = "Mail"

# --- Row heights (as authored after the reshape) -------------------------
$ws.Rows("2").RowHeight = 45
$ws.Rows("3").RowHeight = 30
$ws.Rows("4").RowHeight = 30

# --- Column widths ---------------------------------------------------------
$ws.Columns("A").ColumnWidth = 20.666666666666668
$ws.Columns("B").ColumnWidth = 21.0
$ws.Columns("C").ColumnWidth = 26.833333333333332
$ws.Columns("D").ColumnWidth = 30.666666666666668
$ws.Columns("E").ColumnWidth = 10.666666666666666
$ws.Columns("F").ColumnWidth = 14.666666666666666

# --- Selection / view tidy-up --------------------------------------------
$ws.Range("D11").Select()
